# span_060 commit "span 60 done in 3 jobs": the Q-vs-laser-power sweep for the
# 60 km span was captured in three simulation batches. The first batch (Length of
# Segment 1 = 0 and 6 km, rows 2-25) was already complete. The second batch left
# partial sweeps for 12/18/24 km (rows 26-49 before the edit, each group missing
# the tail/head of its -7..4 dBm sweep). This change fills in those missing A/Q
# pairs and appends the third batch's full sweeps for 30/36/42/48/54/60 km,
# growing the sheet from 119 to 133 rows (dimension A1:I119 -> A1:I133).
# Rows 1-33 are untouched; rows 34-133 are (re)written below in the canonical
# -7,-6,-5,-4,-3,-2,-1,0,1,2,3,4 order for each Length-of-Segment-1 group.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$aVals = @(1,2,3,4,-7,-6,-5,-4,-3,-2,-1,0,1,2,3,4,-7,-6,-5,-4,-3,-2,-1,0,1,2,3,4,-7,-6,-5,-4,-3,-2,-1,0,1,2,3,4,-7,-6,-5,-4,-3,-2,-1,0,1,2,3,4,-7,-6,-5,-4,-3,-2,-1,0,1,2,3,4,-7,-6,-5,-4,-3,-2,-1,0,1,2,3,4,-7,-6,-5,-4,-3,-2,-1,0,1,2,3,4,-7,-6,-5,-4,-3,-2,-1,0,1,2,3,4)
$bVals = @(12,12,12,12,18,18,18,18,18,18,18,18,18,18,18,18,24,24,24,24,24,24,24,24,24,24,24,24,30,30,30,30,30,30,30,30,30,30,30,30,36,36,36,36,36,36,36,36,36,36,36,36,42,42,42,42,42,42,42,42,42,42,42,42,48,48,48,48,48,48,48,48,48,48,48,48,54,54,54,54,54,54,54,54,54,54,54,54,60,60,60,60,60,60,60,60,60,60,60,60)
$cVals = @(4.2183,3.2191,2.2645,0.08040599999999999,5.4061,6.1348,6.7679,7.3124,7.6755,7.9194,7.4727,7.114,5.7714,3.9522,2.7962,0.98326,5.2224,5.9398,6.5637,7.2437,7.6557,7.7875,7.7265,7.3157,6.7368,5.1046,3.9382,2.1191,5.176,5.7738,6.5183,7.0448,7.5354,7.816,7.7553,7.5713,6.9593,5.9777,4.9941,2.6788,4.988,5.6586,6.1847,6.8691,7.4006,7.7265,7.9058,7.6066,7.0544,5.4424,5.2328,2.2253,4.9102,5.5181,6.1269,6.7209,7.2649,7.5838,7.6901,7.6874,7.2681,5.8734,5.577999999999999,3.5557,4.7554,5.3068,6.035,6.5754,7.0943,7.499,7.6169,7.5038,7.229,5.6217,5.2794,3.8022,4.5935,5.2011,5.8226,6.5321,6.9319,7.3069,7.5074,7.6323,7.3211,6.7061,5.8512,3.7773,4.5165,5.0926,5.6608,6.2076,6.7209,7.1917,7.4763,7.511,7.0747,6.6929,5.8297,4.3487)

$constD = 60
$constE = 0
$constF = 250
$constG = 112
$constH = 0.16
$constI = 0.158

$startRow = 34
for ($i = 0; $i -lt $aVals.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $aVals[$i]
    $ws.Cells.Item($r, 2).Value = $bVals[$i]
    $ws.Cells.Item($r, 3).Value = $cVals[$i]
    $ws.Cells.Item($r, 4).Value = $constD
    $ws.Cells.Item($r, 5).Value = $constE
    $ws.Cells.Item($r, 6).Value = $constF
    $ws.Cells.Item($r, 7).Value = $constG
    $ws.Cells.Item($r, 8).Value = $constH
    $ws.Cells.Item($r, 9).Value = $constI
}

$lastRow = $startRow + $aVals.Length - 1
Write-Output "wrote rows $startRow..$lastRow (dimension now A1:I$lastRow)"
